$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 gets its own (non-shared) formula
$ws.Range("D2").Formula = '=CONCATENATE(A2,"=",B2,"~",C2)'

# D3:D20 entered together, forming a shared formula group
$ws.Range("D3:D20").Formula = '=CONCATENATE(A3,"=",B3,"~",C3)'

# Set column D width to match bestFit width of other imported data
$ws.Columns.Item(4).ColumnWidth = 25

# Update the selection to match the new active range
$ws.Range("D2:D20").Select()
